# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# Updates specific profit-calculation cells (H, I, J, K, L, M, N) across several sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 200
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 200
$ws.Range("N12").Value = -540
$ws.Range("M12").ClearContents()

$ws.Range("H48").Value = 1183.3334
$ws.Range("J48").Value = 1525
$ws.Range("L48").Value = 4575
$ws.Range("N48").Value = -5159

$ws.Range("H56").Value = 1183.3334
$ws.Range("J56").Value = 1525
$ws.Range("L56").Value = 4575
$ws.Range("N56").Value = -5643

$ws.Range("H125").Value = 1938.4
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1938.4
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 17445.6
$ws.Range("N125").Value = -22365.6
$ws.Range("M125").ClearContents()

$ws.Range("H129").Value = 999.06665
$ws.Range("J129").Value = 1147.7826
$ws.Range("L129").Value = 3443.3478
$ws.Range("N129").Value = -13443.3478

$ws.Range("H132").Value = 2505.721
$ws.Range("I132").Value = 2073.675
$ws.Range("J132").Value = 8266.333000000001
$ws.Range("K132").Value = 6221.025000000001
$ws.Range("L132").Value = 24798.999
$ws.Range("M132").Value = -3691.025000000001
$ws.Range("N132").Value = -29858.999

$ws.Range("H135").Value = 790.4318
$ws.Range("I135").Value = 348.83783
$ws.Range("K135").Value = 3139.54047
$ws.Range("M135").Value = -604.5404699999999

$ws.Range("H138").Value = 2334.963
$ws.Range("I138").Value = 2116.318
$ws.Range("J138").Value = 3297
$ws.Range("K138").Value = 6348.954000000001
$ws.Range("L138").Value = 9891
$ws.Range("M138").Value = -1208.954000000001
$ws.Range("N138").Value = -20171

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 539078.8
$ws.Range("I32").Value = 630661.6
$ws.Range("J32").Value = 17056.8
$ws.Range("K32").Value = 630661.6
$ws.Range("L32").Value = 17056.8
$ws.Range("M32").Value = -630374.6
$ws.Range("N32").Value = -17630.8

$ws.Range("H34").Value = 70028
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws.Range("H61").Value = 2044.1136
$ws.Range("I61").Value = 1656.6945
$ws.Range("K61").Value = 1656.6945
$ws.Range("M61").Value = -1444.6945

$ws.Range("H132").Value = 2900.4062
$ws.Range("I132").Value = 1841.6086
$ws.Range("J132").Value = 5606.222
$ws.Range("K132").Value = 5524.825800000001
$ws.Range("L132").Value = 16818.666
$ws.Range("M132").Value = -2994.825800000001
$ws.Range("N132").Value = -21878.666

$ws.Range("H136").Value = 2044.1136
$ws.Range("I136").Value = 1656.6945
$ws.Range("K136").Value = 4970.083500000001
$ws.Range("M136").Value = -2420.083500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 34722.168
$ws.Range("J88").Value = 34722.168
$ws.Range("L88").Value = 34722.168
$ws.Range("N88").Value = -35534.168

$ws.Range("H91").Value = 34722.168
$ws.Range("J91").Value = 34722.168
$ws.Range("L91").Value = 34722.168
$ws.Range("N91").Value = -37530.168

$ws.Range("H134").Value = 2676.3547
$ws.Range("I134").Value = 2270.15
$ws.Range("K134").Value = 6810.450000000001
$ws.Range("M134").Value = -4275.450000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5207.8623
$ws.Range("I31").Value = 1284.4688
$ws.Range("J31").Value = 10036.654
$ws.Range("K31").Value = 1284.4688
$ws.Range("L31").Value = 10036.654
$ws.Range("M31").Value = -989.4688000000001
$ws.Range("N31").Value = -10626.654

$ws.Range("H34").Value = 5207.8623
$ws.Range("I34").Value = 1284.4688
$ws.Range("J34").Value = 10036.654
$ws.Range("K34").Value = 1284.4688
$ws.Range("L34").Value = 10036.654
$ws.Range("M34").Value = -1082.4688
$ws.Range("N34").Value = -10440.654

$ws.Range("H58").Value = 1648.4286
$ws.Range("I58").Value = 1325.2
$ws.Range("J58").Value = 2021.3846
$ws.Range("K58").Value = 1325.2
$ws.Range("L58").Value = 2021.3846
$ws.Range("M58").Value = -1122.2
$ws.Range("N58").Value = -2427.3846

$ws.Range("H132").Value = 3402759.5
$ws.Range("I132").Value = 1072.5
$ws.Range("K132").Value = 3217.5
$ws.Range("M132").Value = -687.5

$ws.Range("H134").Value = 3674.077
$ws.Range("I134").Value = 4270.5
$ws.Range("J134").Value = 1686
$ws.Range("K134").Value = 12811.5
$ws.Range("L134").Value = 5058
$ws.Range("M134").Value = -10276.5
$ws.Range("N134").Value = -10128

$ws.Range("H136").Value = 1648.4286
$ws.Range("I136").Value = 1325.2
$ws.Range("J136").Value = 2021.3846
$ws.Range("K136").Value = 3975.6
$ws.Range("L136").Value = 6064.1538
$ws.Range("M136").Value = -1425.6
$ws.Range("N136").Value = -11164.1538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 768.6842
$ws.Range("I5").Value = 564.46155
$ws.Range("J5").Value = 1211.1666
$ws.Range("K5").Value = 1693.38465
$ws.Range("L5").Value = 3633.4998
$ws.Range("M5").Value = -1581.38465
$ws.Range("N5").Value = -3857.4998

$ws.Range("H81").Value = 5203
$ws.Range("J81").Value = 9499.5
$ws.Range("L81").Value = 28498.5
$ws.Range("N81").Value = -30744.5

$ws.Range("H84").Value = 5203
$ws.Range("J84").Value = 9499.5
$ws.Range("L84").Value = 85495.5
$ws.Range("N84").Value = -96727.5

$ws.Range("H113").Value = 695.73914
$ws.Range("I113").Value = 501.97437
$ws.Range("J113").Value = 947.63336
$ws.Range("K113").Value = 1505.92311
$ws.Range("L113").Value = 2842.90008
$ws.Range("M113").Value = 664.0768899999998
$ws.Range("N113").Value = -7182.90008

$ws.Range("H132").Value = 2370.6743
$ws.Range("I132").Value = 2417.111
$ws.Range("J132").Value = 2358.3823
$ws.Range("K132").Value = 21753.999
$ws.Range("L132").Value = 21225.4407
$ws.Range("M132").Value = -19223.999
$ws.Range("N132").Value = -26285.4407

$ws.Range("H135").Value = 768.6842
$ws.Range("I135").Value = 564.46155
$ws.Range("J135").Value = 1211.1666
$ws.Range("K135").Value = 5080.15395
$ws.Range("L135").Value = 10900.4994
$ws.Range("M135").Value = -2545.15395
$ws.Range("N135").Value = -15970.4994

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1621.48
$ws.Range("I132").Value = 1185.4474
$ws.Range("J132").Value = 3002.25
$ws.Range("K132").Value = 3556.3422
$ws.Range("L132").Value = 9006.75
$ws.Range("M132").Value = -1026.3422
$ws.Range("N132").Value = -14066.75

$ws.Range("H136").Value = 4903127
$ws.Range("I136").Value = 1004.55554
$ws.Range("J136").Value = 23811314
$ws.Range("K136").Value = 3013.66662
$ws.Range("L136").Value = 71433942
$ws.Range("M136").Value = -463.66662
$ws.Range("N136").Value = -71439042

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 302508.25
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 302508.25
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 302508.25
$ws.Range("N29").Value = -303088.25
$ws.Range("M29").ClearContents()

$ws.Range("H82").Value = 39980
$ws.Range("J82").Value = 39980
$ws.Range("L82").Value = 39980
$ws.Range("N82").Value = -40746

$ws.Range("H85").Value = 39980
$ws.Range("J85").Value = 39980
$ws.Range("L85").Value = 39980
$ws.Range("N85").Value = -42632

$ws.Range("H122").Value = 1894.579
$ws.Range("I122").Value = 1853.2142
$ws.Range("K122").Value = 5559.642599999999
$ws.Range("M122").Value = -3109.642599999999

$ws.Range("H136").Value = 1746.6842
$ws.Range("I136").Value = 1266.1555
$ws.Range("J136").Value = 3548.6667
$ws.Range("K136").Value = 3798.4665
$ws.Range("L136").Value = 10646.0001
$ws.Range("M136").Value = -1248.4665
